$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NovemberRaw")

# Header row
$ws.Range("A1").Value = 'Library'
$ws.Range("B1").Value = 'Items owned by this library checked out at this library this month'
$ws.Range("C1").Value = 'Items owned by other libraries checked out at this library this month'
$ws.Range("D1").Value = 'Total circulation this month'

# Library rows
$ws.Range("A2").Value = 'Atchison Public Library'
$ws.Range("B2").Value = 3894
$ws.Range("C2").Value = 1468
$ws.Range("D2").Value = 5362
$ws.Range("A3").Value = 'Baldwin City Public Library'
$ws.Range("B3").Value = 2403
$ws.Range("C3").Value = 450
$ws.Range("D3").Value = 2853
$ws.Range("A4").Value = 'Basehor Community Library'
$ws.Range("B4").Value = 7590
$ws.Range("C4").Value = 1092
$ws.Range("D4").Value = 8682
$ws.Range("A5").Value = 'Bern Community Library'
$ws.Range("B5").Value = 130
$ws.Range("C5").Value = 22
$ws.Range("D5").Value = 152
$ws.Range("A6").Value = 'Bonner Springs City Library'
$ws.Range("B6").Value = 4779
$ws.Range("C6").Value = 1141
$ws.Range("D6").Value = 5920
$ws.Range("A7").Value = 'Burlingame Community Library'
$ws.Range("B7").Value = 420
$ws.Range("C7").Value = 243
$ws.Range("D7").Value = 663
$ws.Range("A8").Value = 'Carbondale City Library'
$ws.Range("B8").Value = 470
$ws.Range("C8").Value = 70
$ws.Range("D8").Value = 540
$ws.Range("A9").Value = 'Centralia Community Library'
$ws.Range("B9").Value = 310
$ws.Range("C9").Value = 70
$ws.Range("D9").Value = 380
$ws.Range("A10").Value = 'Corning City Library'
$ws.Range("B10").Value = 9
$ws.Range("C10").Value = 2
$ws.Range("D10").Value = 11
$ws.Range("A11").Value = 'Digital Content'
$ws.Range("A12").Value = 'Doniphan County Library - Elwood'
$ws.Range("B12").Value = 144
$ws.Range("C12").Value = 41
$ws.Range("D12").Value = 185
$ws.Range("A13").Value = 'Doniphan County Library - Highland'
$ws.Range("B13").Value = 272
$ws.Range("C13").Value = 182
$ws.Range("D13").Value = 454
$ws.Range("A14").Value = 'Doniphan County Library - Troy'
$ws.Range("B14").Value = 446
$ws.Range("C14").Value = 128
$ws.Range("D14").Value = 574
$ws.Range("A15").Value = 'Doniphan County Library - Wathena'
$ws.Range("B15").Value = 419
$ws.Range("C15").Value = 97
$ws.Range("D15").Value = 516
$ws.Range("A16").Value = 'Effingham Community Library'
$ws.Range("B16").Value = 192
$ws.Range("C16").Value = 19
$ws.Range("D16").Value = 211
$ws.Range("A17").Value = 'Eudora Community Library'
$ws.Range("B17").Value = 1499
$ws.Range("C17").Value = 593
$ws.Range("D17").Value = 2092
$ws.Range("A18").Value = 'Everest, Barnes Reading Room'
$ws.Range("B18").Value = 110
$ws.Range("C18").Value = 67
$ws.Range("D18").Value = 177
$ws.Range("A19").Value = 'Hiawatha, Morrill Public Library'
$ws.Range("B19").Value = 1606
$ws.Range("C19").Value = 629
$ws.Range("D19").Value = 2235
$ws.Range("A20").Value = 'Highland Community College'
$ws.Range("B20").Value = 32
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 36
$ws.Range("A21").Value = 'Holton, Beck-Bookman Library'
$ws.Range("B21").Value = 1750
$ws.Range("C21").Value = 405
$ws.Range("D21").Value = 2155
$ws.Range("A22").Value = 'Horton Public Library'
$ws.Range("B22").Value = 90
$ws.Range("C22").Value = 31
$ws.Range("D22").Value = 121
$ws.Range("A23").Value = 'Lansing Community Library'
$ws.Range("B23").Value = 1736
$ws.Range("C23").Value = 766
$ws.Range("D23").Value = 2502
$ws.Range("A24").Value = 'Leavenworth Public Library'
$ws.Range("B24").Value = 8176
$ws.Range("C24").Value = 2023
$ws.Range("D24").Value = 10199
$ws.Range("A25").Value = 'Linwood Community Library'
$ws.Range("B25").Value = 567
$ws.Range("C25").Value = 208
$ws.Range("D25").Value = 775
$ws.Range("A26").Value = 'Louisburg Library'
$ws.Range("A27").Value = 'Lyndon Carnegie Library'
$ws.Range("B27").Value = 411
$ws.Range("C27").Value = 278
$ws.Range("D27").Value = 689
$ws.Range("A28").Value = 'McLouth Public Library'
$ws.Range("B28").Value = 160
$ws.Range("C28").Value = 44
$ws.Range("D28").Value = 204
$ws.Range("A29").Value = 'Meriden-Ozawkie Public Library'
$ws.Range("B29").Value = 1275
$ws.Range("C29").Value = 570
$ws.Range("D29").Value = 1845
$ws.Range("A30").Value = 'Northeast Kansas Library System'
$ws.Range("B30").Value = 10
$ws.Range("C30").Value = 24
$ws.Range("D30").Value = 34
$ws.Range("A31").Value = 'Nortonville Public Library'
$ws.Range("B31").Value = 260
$ws.Range("C31").Value = 86
$ws.Range("D31").Value = 346
$ws.Range("A32").Value = 'Osage City Library'
$ws.Range("B32").Value = 1421
$ws.Range("C32").Value = 419
$ws.Range("D32").Value = 1840
$ws.Range("A33").Value = 'Osawatomie Public Library'
$ws.Range("B33").Value = 782
$ws.Range("C33").Value = 344
$ws.Range("D33").Value = 1126
$ws.Range("A34").Value = 'Oskaloosa Public Library'
$ws.Range("B34").Value = 473
$ws.Range("C34").Value = 216
$ws.Range("D34").Value = 689
$ws.Range("A35").Value = 'Ottawa Library'
$ws.Range("B35").Value = 5953
$ws.Range("C35").Value = 853
$ws.Range("D35").Value = 6806
$ws.Range("A36").Value = 'Overbrook Public Library'
$ws.Range("B36").Value = 921
$ws.Range("C36").Value = 192
$ws.Range("D36").Value = 1113
$ws.Range("A37").Value = 'Paola Free Library'
$ws.Range("B37").Value = 3133
$ws.Range("C37").Value = 435
$ws.Range("D37").Value = 3568
$ws.Range("A38").Value = 'Perry-Lecompton Community Library'
$ws.Range("B38").Value = 80
$ws.Range("C38").Value = 19
$ws.Range("D38").Value = 99
$ws.Range("A39").Value = 'Pomona Community Library'
$ws.Range("B39").Value = 104
$ws.Range("C39").Value = 74
$ws.Range("D39").Value = 178
$ws.Range("A40").Value = 'Prairie Hills Schools - Axtell Public School'
$ws.Range("B40").Value = 488
$ws.Range("C40").Value = 35
$ws.Range("D40").Value = 523
$ws.Range("A41").Value = 'Prairie Hills Schools - Sabetha Elementary School'
$ws.Range("B41").Value = 2079
$ws.Range("C41").Value = 86
$ws.Range("D41").Value = 2165
$ws.Range("A42").Value = 'Prairie Hills Schools - Sabetha High School'
$ws.Range("B42").Value = 30
$ws.Range("C42").Value = 2
$ws.Range("D42").Value = 32
$ws.Range("A43").Value = 'Prairie Hills Schools - Sabetha Middle School'
$ws.Range("B43").Value = 168
$ws.Range("C43").Value = 15
$ws.Range("D43").Value = 183
$ws.Range("A44").Value = 'Prairie Hills Schools - Wetmore Academic Center (Permanently closed)'
$ws.Range("A45").Value = 'Richmond Public Library'
$ws.Range("B45").Value = 415
$ws.Range("C45").Value = 77
$ws.Range("D45").Value = 492
$ws.Range("A46").Value = 'Rossville Community Library'
$ws.Range("B46").Value = 1271
$ws.Range("C46").Value = 507
$ws.Range("D46").Value = 1778
$ws.Range("A47").Value = 'Sabetha, Mary Cotton Library'
$ws.Range("B47").Value = 2903
$ws.Range("C47").Value = 930
$ws.Range("D47").Value = 3833
$ws.Range("A48").Value = 'Seneca Free Library'
$ws.Range("B48").Value = 1404
$ws.Range("C48").Value = 229
$ws.Range("D48").Value = 1633
$ws.Range("A49").Value = 'Silver Lake Library'
$ws.Range("B49").Value = 1025
$ws.Range("C49").Value = 457
$ws.Range("D49").Value = 1482
$ws.Range("A50").Value = 'Tonganoxie Public Library'
$ws.Range("B50").Value = 3003
$ws.Range("C50").Value = 1043
$ws.Range("D50").Value = 4046
$ws.Range("A51").Value = 'Valley Falls, Delaware Township Library'
$ws.Range("B51").Value = 377
$ws.Range("C51").Value = 163
$ws.Range("D51").Value = 540
$ws.Range("A52").Value = 'Wellsville City Library'
$ws.Range("B52").Value = 1095
$ws.Range("C52").Value = 458
$ws.Range("D52").Value = 1553
$ws.Range("A53").Value = 'Wetmore Public Library'
$ws.Range("B53").Value = 163
$ws.Range("C53").Value = 100
$ws.Range("D53").Value = 263
$ws.Range("A54").Value = 'Williamsburg Community Library'
$ws.Range("B54").Value = 399
$ws.Range("C54").Value = 3
$ws.Range("D54").Value = 402
$ws.Range("A55").Value = 'Winchester Public Library'
$ws.Range("B55").Value = 365
$ws.Range("C55").Value = 340
$ws.Range("D55").Value = 705
